# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "25.855.16"
    "E2"  = "  -0.36%  "
    "D3"  = "1.585.46"
    "E3"  = "  -2.23%  "
    "E4"  = "  +0.00%  "
    "D5"  = "209.93"
    "E5"  = "  -1.24%  "
    "E6"  = "  -0.01%  "
    "D7"  = "0.477"
    "E7"  = "  -3.81%  "
    "E8"  = "  -0.89%  "
    "E9"  = "  -0.50%  "
    "D10" = "18.04"
    "E10" = "  -2.06%  "
    "E11" = "  +0.03%  "
    "D12" = "1.805.55"
    "E12" = "  -2.24%  "
    "D13" = "1.584.96"
    "E13" = "  -2.14%  "
    "D14" = "4.03"
    "E14" = "  -2.81%  "
    "E15" = "  -3.00%  "
    "D16" = "25.827.15"
    "E16" = "  -0.53%  "
    "D17" = "0.0₃0724"
    "E17" = "  -2.12%  "
    "D18" = "59.71"
    "E18" = "  -3.31%  "
    "E19" = "  -0.03%  "
    "D20" = "191.55"
    "E20" = "  -0.38%  "
    "E21" = "  -1.78%  "
    "E22" = "  -1.82%  "
    "E23" = "  -1.37%  "
    "E24" = "  -1.11%  "
    "D25" = "142.02"
    "E25" = "  -1.35%  "
    "E26" = "  +0.04%  "
    "E27" = "  -0.39%  "
    "D28" = "15.08"
    "E28" = "  -1.02%  "
    "E29" = "  -2.96%  "
    "E30" = "  -5.41%  "
    "E31" = "  -1.59%  "
    "E33" = "  -2.58%  "
    "E34" = "  -0.20%  "
    "D35" = "2.36"
    "E35" = "  -2.33%  "
    "D36" = "1.099.76"
    "E36" = "  -2.50%  "
    "E37" = "  +0.04%  "
    "E38" = "  -2.02%  "
    "D39" = "0.501"
    "E39" = "  -2.93%  "
    "E40" = "  -2.17%  "
    "E41" = "  +8.07%  "
    "D42" = "0.776"
    "E42" = "  -8.17%  "
    "E43" = "  +1.40%  "
    "D44" = "93.90"
    "E44" = "  -3.98%  "
    "E45" = "  -2.19%  "
    "D46" = "0.0₆0111"
    "E46" = "  -1.31%  "
    "E47" = "  -0.51%  "
    "E49" = "  -1.62%  "
    "D50" = "0.408"
    "E50" = "  -0.62%  "
    "E51" = "  +0.02%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.Value = "'" + $updates[$ref]
    $cell.Style = "Normal"
}
